$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("KDR<0.001")

$shp = $ws1.Shapes.Item(1)
# offset of the shape's left edge from the left edge of its anchor column (col G, 1-indexed 7)
$colOffset = $shp.Left - $ws1.Cells.Item(1, 7).Left
$shpTop = $shp.Top
$shpWidth = $shp.Width
$shpHeight = $shp.Height

$ws1.Columns("E:E").Delete()

$ws1.Range("C1").Value = "tau"
$ws1.Range("B1").Value = "n"
$ws1.Range("A1").Value = "pch"
$ws1.Range("D1").Value = "R"

$shp2 = $ws1.Shapes.Item(1)
$shp2.Left = $ws1.Cells.Item(1, 6).Left + $colOffset
$shp2.Top = $shpTop
$shp2.Width = $shpWidth
$shp2.Height = $shpHeight

$ws1.Activate() | Out-Null
$ws1.Range("L26").Select() | Out-Null
